$wb = $excel.ActiveWorkbook

# Source sheet to clone the layout/styles from (groupDirNone), which already
# has the exact same dimensions/column widths/styles as the two new sheets.
$src = $wb.Worksheets.Item("groupDirNone")

# --- Add "indexVar" sheet (copied after the last existing sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $lastSheet)
$indexVarSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$indexVarSheet.Name = "indexVar"
$indexVarSheet.Range("A3").Value = '${index + 1}. ${divisionsList.teams.city}?@indexVar=index'
[void]$indexVarSheet.Range("A3").Select()

# --- Add "limit" sheet (copied after "indexVar") ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $lastSheet2)
$limitSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$limitSheet.Name = "limit"
$limitSheet.Range("A3").Value = '${divisionsList.teams.city}?@limit=3'
